$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mobile numbers in column D (rows 2-5) to unique values per row
$ws.Range("D2").Value = "416-858-7781"
$ws.Range("D3").Value = "416-858-7782"
$ws.Range("D4").Value = "416-858-7783"
$ws.Range("D5").Value = "416-858-7784"

# Drop the stale hyperlinks so they can be rebuilt against the new (shrunk) ranges
$ws.Hyperlinks.Delete()

# Remove row 6 (the long-password / saritha756 test row) entirely
$ws.Rows("6:6").Delete()

# Give B5 (the numeric mobile-style cell) a left-aligned look
$ws.Range("B5").HorizontalAlignment = -4131

# Rebuild the hyperlinks against the surviving rows (2-5)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:saritha750@yahoo.ca")
$ws.Hyperlinks.Add($ws.Range("C3:C5"), "mailto:saritha750@yahoo.ca", "", "", "saritha750@yahoo.ca")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:saritha751@yahoo.ca")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:saritha753@yahoo.ca")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:saritha755@yahoo.ca")

# Widen columns A-C to match the resized layout
$ws.Columns.Item(1).ColumnWidth = 15.42578125
$ws.Columns.Item(2).ColumnWidth = 21.140625
$ws.Columns.Item(3).ColumnWidth = 26.140625

# Move the active selection to C10, matching the saved view state
[void]$ws.Range("C10").Select()
